# Commit: "Updated the input files in sixteen_tests to have
# strain_log2_expression instead of just strain and then ran the files and
# saved the outputs in sixteen_tests_output"
#
# Rename the two per-strain expression sheets to make their contents
# ("log2 expression") explicit in the tab name.
$wb = $excel.ActiveWorkbook

$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

# Remember which sheet/tab was active before we touch anything, so the
# workbook's active tab is unchanged once we're done.
$originalActiveSheet = $wb.ActiveSheet.Name

# The re-saved workbook shows the selection on "wt_log2_expression" moved
# from N1:N5 to B45 (e.g. from scrolling/clicking around in Excel while
# inspecting the renamed sheet).
$wsWt.Activate()
$wsWt.Range("B45").Select()

# Restore the originally active sheet/tab.
$wb.Worksheets.Item($originalActiveSheet).Activate()
